$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 387.16666
$ws.Range("I9").Value = 204.88889
$ws.Range("J9").Value = 934
$ws.Range("K9").Value = 204.88889
$ws.Range("L9").Value = 934
$ws.Range("M9").Value = -35.88889
$ws.Range("N9").Value = -1272

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 22387.555
$ws.Range("J43").Value = 26498
$ws.Range("L43").Value = 26498
$ws.Range("N43").Value = -26636

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 315
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 11108.333
$ws.Range("I70").Value = 18000
$ws.Range("J70").Value = 9139.286
$ws.Range("K70").Value = 54000
$ws.Range("L70").Value = 27417.858
$ws.Range("M70").Value = -53730
$ws.Range("N70").Value = -27957.858

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 11108.333
$ws.Range("I73").Value = 18000
$ws.Range("J73").Value = 9139.286
$ws.Range("K73").Value = 54000
$ws.Range("L73").Value = 27417.858
$ws.Range("M73").Value = -53064
$ws.Range("N73").Value = -29289.858

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 33348.6
$ws.Range("J87").Value = 33348.6
$ws.Range("L87").Value = 33348.6
$ws.Range("N87").Value = -35844.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 33348.6
$ws.Range("J90").Value = 33348.6
$ws.Range("L90").Value = 100045.8
$ws.Range("N90").Value = -112525.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 31251410
$ws.Range("I103").Value = 1181
$ws.Range("K103").Value = 3543
$ws.Range("M103").Value = -2957

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2168.7083
$ws.Range("J112").Value = 2206.739
$ws.Range("L112").Value = 6620.217000000001
$ws.Range("N112").Value = -8836.217000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1094.4445
$ws.Range("I125").Value = 981.25
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 8831.25
$ws.Range("L125").Value = 18000
$ws.Range("M125").Value = -6371.25
$ws.Range("N125").Value = -22920

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4932.227
$ws.Range("J138").Value = 5071.702
$ws.Range("L138").Value = 15215.106
$ws.Range("N138").Value = -25495.106

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10761.172
$ws.Range("I32").Value = 8367.552
$ws.Range("J32").Value = 33899.5
$ws.Range("K32").Value = 8367.552
$ws.Range("L32").Value = 33899.5
$ws.Range("M32").Value = -8080.552
$ws.Range("N32").Value = -34473.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3164.6667
$ws.Range("I61").Value = 2248.5
$ws.Range("J61").Value = 4997
$ws.Range("K61").Value = 2248.5
$ws.Range("L61").Value = 4997
$ws.Range("M61").Value = -2036.5
$ws.Range("N61").Value = -5421

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 105610.4
$ws.Range("J118").Value = 105610.4
$ws.Range("L118").Value = 105610.4
$ws.Range("N118").Value = -108924.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3202.276
$ws.Range("I132").Value = 3191.6428
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 9574.928400000001
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -7044.928400000001
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3164.6667
$ws.Range("I136").Value = 2248.5
$ws.Range("J136").Value = 4997
$ws.Range("K136").Value = 6745.5
$ws.Range("L136").Value = 14991
$ws.Range("M136").Value = -4195.5
$ws.Range("N136").Value = -20091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 98180.09
$ws.Range("I134").Value = 104858.07
$ws.Range("J134").Value = 1349.5
$ws.Range("K134").Value = 314574.21
$ws.Range("L134").Value = 4048.5
$ws.Range("M134").Value = -312039.21
$ws.Range("N134").Value = -9118.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28080.293
$ws.Range("I31").Value = 58174
$ws.Range("J31").Value = 4528.696
$ws.Range("K31").Value = 58174
$ws.Range("L31").Value = 4528.696
$ws.Range("M31").Value = -57879
$ws.Range("N31").Value = -5118.696

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 28080.293
$ws.Range("I34").Value = 58174
$ws.Range("J34").Value = 4528.696
$ws.Range("K34").Value = 58174
$ws.Range("L34").Value = 4528.696
$ws.Range("M34").Value = -57972
$ws.Range("N34").Value = -4932.696

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2938.1667
$ws.Range("J105").Value = 4000
$ws.Range("L105").Value = 4000
$ws.Range("N105").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3084
$ws.Range("I132").Value = 2823.9023
$ws.Range("K132").Value = 8471.706900000001
$ws.Range("M132").Value = -5941.706900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1612.4445
$ws.Range("I5").Value = 706
$ws.Range("J5").Value = 2337.6
$ws.Range("K5").Value = 2118
$ws.Range("L5").Value = 7012.799999999999
$ws.Range("M5").Value = -2006
$ws.Range("N5").Value = -7236.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1600.3334
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 6000
$ws.Range("M68").Value = -5189

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1600.3334
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 18000
$ws.Range("M71").Value = -13944

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 103771.4
$ws.Range("I94").Value = 334079.66
$ws.Range("J94").Value = 5067.857
$ws.Range("K94").Value = 1002238.98
$ws.Range("L94").Value = 15203.571
$ws.Range("M94").Value = -1001562.98
$ws.Range("N94").Value = -16555.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 7306.9165
$ws.Range("I100").Value = 3604
$ws.Range("J100").Value = 9951.857
$ws.Range("K100").Value = 10812
$ws.Range("L100").Value = 29855.571
$ws.Range("M100").Value = -10001
$ws.Range("N100").Value = -31477.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 5817.8
$ws.Range("I104").Value = 2550
$ws.Range("J104").Value = 7996.3335
$ws.Range("K104").Value = 7650
$ws.Range("L104").Value = 23989.0005
$ws.Range("M104").Value = -5029
$ws.Range("N104").Value = -29231.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4356.9
$ws.Range("I113").Value = 2750
$ws.Range("J113").Value = 4535.4443
$ws.Range("K113").Value = 8250
$ws.Range("L113").Value = 13606.3329
$ws.Range("M113").Value = -6080
$ws.Range("N113").Value = -17946.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1581.6
$ws.Range("I122").Value = 1169.3334
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 10524.0006
$ws.Range("L122").Value = 19800
$ws.Range("M122").Value = -8074.000599999999
$ws.Range("N122").Value = -24700

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1612.4445
$ws.Range("I135").Value = 706
$ws.Range("J135").Value = 2337.6
$ws.Range("K135").Value = 6354
$ws.Range("L135").Value = 21038.4
$ws.Range("M135").Value = -3819
$ws.Range("N135").Value = -26108.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1049.4667
$ws.Range("I2").Value = 611.1429000000001
$ws.Range("K2").Value = 611.1429000000001
$ws.Range("M2").Value = -498.1429000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2000000
$ws.Range("I14").Value = 2000000
$ws.Range("K14").Value = 2000000
$ws.Range("M14").Value = -1999832

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3059
$ws.Range("J102").Value = 3450
$ws.Range("L102").Value = 3450
$ws.Range("N102").Value = -6694

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 45066.25
$ws.Range("J106").Value = 45066.25
$ws.Range("L106").Value = 45066.25
$ws.Range("N106").Value = -47590.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 134996.17
$ws.Range("J133").Value = 134996.17
$ws.Range("L133").Value = 134996.17
$ws.Range("N133").Value = -145116.17

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 122665.14
$ws.Range("I61").Value = 87553.836
$ws.Range("K61").Value = 87553.836
$ws.Range("M61").Value = -87351.836

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 16637.934
$ws.Range("J93").Value = 42301.223
$ws.Range("L93").Value = 42301.223
$ws.Range("N93").Value = -44797.223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 122665.14
$ws.Range("I113").Value = 87553.836
$ws.Range("K113").Value = 87553.836
$ws.Range("M113").Value = -85383.836

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6497.1665
$ws.Range("I132").Value = 5747
$ws.Range("J132").Value = 7997.5
$ws.Range("K132").Value = 17241
$ws.Range("L132").Value = 23992.5
$ws.Range("M132").Value = -14711
$ws.Range("N132").Value = -29052.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5083.2354
$ws.Range("I132").Value = 5213.4375
$ws.Range("K132").Value = 15640.3125
$ws.Range("M132").Value = -13110.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1747
$ws.Range("I136").Value = 1496.6666
$ws.Range("J136").Value = 1997.3334
$ws.Range("K136").Value = 4489.9998
$ws.Range("L136").Value = 5992.0002
$ws.Range("M136").Value = -1939.9998
$ws.Range("N136").Value = -11092.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 109455.875
$ws.Range("I139").Value = 650
$ws.Range("J139").Value = 124999.57
$ws.Range("K139").Value = 650
$ws.Range("L139").Value = 124999.57
$ws.Range("M139").Value = 4490
$ws.Range("N139").Value = -135279.57
